$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell -> new value for column F ("dSF"), per diff
$updates = @{
    "F5"  = -6
    "F11" = -7
    "F14" = -9
    "F15" = -6
    "F26" = 4
    "F28" = -1
    "F33" = 12
    "F35" = -2
    "F37" = 2
    "F39" = -1
    "F40" = 0
    "F42" = -3
    "F45" = 2
    "F47" = 2
    "F53" = 4
    "F59" = 0
    "F60" = -5
    "F62" = 3
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
